$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transit_simple")

# New "Mode Simple" (column B) values for rows 2-34, reflecting the
# re-categorization described in the commit (household survey summary update).
$values = @(
    "Walk",    # row 2  - Walk, jog, or wheelchair
    "Bike",    # row 3  - Bicycle or e-bike
    "Drive",   # row 4  - Household vehicle 1
    "Drive",   # row 5  - Household vehicle 2
    "Drive",   # row 6  - Household vehicle 3
    "Drive",   # row 7  - Household vehicle 4
    "Drive",   # row 8  - Household vehicle 5
    "Drive",   # row 9  - Household vehicle 6
    "Drive",   # row 10 - Household vehicle 7
    "Drive",   # row 11 - Household vehicle 8
    "Drive",   # row 12 - Household vehicle 9
    "Drive",   # row 13 - Household vehicle 10
    "Drive",   # row 14 - Other household vehicle
    "Drive",   # row 15 - Rental car
    "Drive",   # row 16 - Carshare vehicle (e.g., Zipcar, Car2Go, RelayRides, etc.)
    "Transit", # row 17 - Vanpool vehicle
    "Drive",   # row 18 - Other non-household vehicle
    "Transit", # row 19 - Bus (public transit)
    "Transit", # row 20 - School bus
    "Transit", # row 21 - Private bus or shuttle
    "Transit", # row 22 - Paratransit
    "Transit", # row 23 - Other bus (rMove only)
    "Other",   # row 24 - Airplane or helicopter
    "Transit", # row 25 - Ferry or water taxi
    "Drive",   # row 26 - Car from work
    "Drive",   # row 27 - Friend/colleague's car
    "Other",   # row 28 - Taxi (e.g., Yello Cab)
    "Other",   # row 29 - Other hired service (e.g., Lyft, Uber)
    "Transit", # row 30 - Commuter rail (Sounder, Amtrak)
    "Transit", # row 31 - Other rail (e.g., streetcar)
    "Other",   # row 32 - Other motorcycle/moped/scooter
    "Transit", # row 33 - Urban rail (e.g., Link light rail, monorail)
    "Other"    # row 34 - Other mode (e.g., skateboard, kayak, motorhome, etc.)
)

# Write row 3 (Bike) ahead of row 2 (Walk) so that new shared-string
# entries are created in the same order ("Bike" before "Walk") as in
# the target workbook, then fill in the remaining rows in order.
$ws.Cells.Item(3, 2).Value = $values[1]
$ws.Cells.Item(2, 2).Value = $values[0]

for ($i = 2; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Update the selected cell to match the recorded view state.
$ws.Range("E17").Select()
